# Update the "Metadata" sheet: URL, Version, Date and Publisher moved from the
# Alvearie/IBM project to the LinuxForHealth project (new release 8.0.0).
$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/copay-exclusion"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# Update the "Elements" sheet: the ele-1/ext-1 invariant text was listed on the
# base "Extension" row's Constraint(s) column (AI2); it only really applies to
# the "Extension.extension" row (AI4, which already carries it), so clear it
# from the base row.
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("AI2").Value = ""
